# "Now view log history" -- add two more sign-in/out log sheets ("2" and "3")
# after the existing "0" and "1" sheets, continuing the same 6-column layout,
# and return the active tab to sheet "0".

$wb = $excel.ActiveWorkbook

$headers = @("name", "to", "from", "tout", "tin", "date")

# Writes $val into $cell as literal text, even when $val looks like a date
# (e.g. "12/3/2019"), mirroring how this log was actually populated (plain
# strings from a script) instead of Excel's normal type-inferring cell entry.
function Set-LogValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Write-HeaderRow($ws) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
    }
}

# --- New sheet "2": Varun signs out with Mr. T / Crowe, 01:50PM-02:00PM, 12/3/2019
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "2"
Write-HeaderRow $ws3
Set-LogValue $ws3.Cells.Item(2, 1) "Varun"
Set-LogValue $ws3.Cells.Item(2, 2) "Crowe"
Set-LogValue $ws3.Cells.Item(2, 3) "Mr. T"
Set-LogValue $ws3.Cells.Item(2, 4) "01:50PM"
Set-LogValue $ws3.Cells.Item(2, 5) "02:00PM"
Set-LogValue $ws3.Cells.Item(2, 6) "12/3/2019"

# --- New sheet "3": fresh sheet, header row only, no entries yet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "3"
Write-HeaderRow $ws4

# Reset the old selection left on sheet "1" back to A1, then make sheet "0"
# the active tab again (it was "1" before this edit).
$wb.Worksheets.Item("1").Range("A1").Select()
$wb.Worksheets.Item("0").Activate()
$wb.Worksheets.Item("0").Range("A1").Select()
